# Regenerate the experiment task-order workbook: each task sheet gets a
# freshly "generated" random order (new CSV filenames) and the sheet tabs
# end up re-ordered (NB, TOL, GNG, RS, vSAT) with new numeric-suffixed names.
#
# NOTE: worksheet object handles in this engine are position-based, so a
# previously-captured reference can silently start pointing at a different
# sheet after any Move() call. To stay safe we always re-fetch sheets by
# name (via $wb.Worksheets.Item(...)) immediately before using them.

$wb = $excel.ActiveWorkbook

# --- Re-order the tabs into: NB, TOL, GNG, RS, vSAT ---
$moveOrder = @(
    "NB_TO-16512555342232397",
    "TOL_TO-16512555342862492",
    "GNG_TO-16512555322461243",
    "RS_TO-16512555342282357",
    "vSAT_TO-1651255534365279"
)
foreach ($nm in $moveOrder) {
    $wb.Worksheets.Item($nm).Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
}

# --- Rename tabs to the newly generated identifiers ---
$wb.Worksheets.Item("NB_TO-16512555342232397").Name   = "NB_TO-16515889667478743"
$wb.Worksheets.Item("TOL_TO-16512555342862492").Name  = "TOL_TO-16515889667974896"
$wb.Worksheets.Item("GNG_TO-16512555322461243").Name  = "GNG_TO-16515889668327491"
$wb.Worksheets.Item("RS_TO-16512555342282357").Name   = "RS_TO-16515889668348305"
$wb.Worksheets.Item("vSAT_TO-1651255534365279").Name  = "vSAT_TO-1651588966908666"

# --- NB task order (rows 2-10, column A index 0-8, column B file name) ---
$nbFiles = @(
    "OB-16515889656128066.csv",
    "TB-16515889662042787.csv",
    "ZB-match_4-16515889642047777.csv",
    "TB-16515889667293346.csv",
    "OB-1651588964895019.csv",
    "ZB-match_9-16515889645777974.csv",
    "TB-16515889658970485.csv",
    "ZB-match_4-16515889643906515.csv",
    "OB-1651588965672083.csv"
)
$sheetNB = $wb.Worksheets.Item("NB_TO-16515889667478743")
for ($i = 0; $i -lt $nbFiles.Length; $i++) {
    $row = $i + 2
    $sheetNB.Cells.Item($row, 2).Value = $nbFiles[$i]
}

# --- TOL task order (rows 2-7, column A index 0-5, column B file name) ---
$tolFiles = @(
    "MM_stims-16515889667638173.csv",
    "ZM_stims-16515889667510111.csv",
    "MM_stims-16515889667811732.csv",
    "ZM_stims-16515889667648628.csv",
    "MM_stims-16515889667965233.csv",
    "ZM_stims-1651588966782172.csv"
)
$sheetTOL = $wb.Worksheets.Item("TOL_TO-16515889667974896")
for ($i = 0; $i -lt $tolFiles.Length; $i++) {
    $row = $i + 2
    $sheetTOL.Cells.Item($row, 2).Value = $tolFiles[$i]
}

# --- GNG task order (rows 2-5, column A index 0-3, column B file name) ---
$gngFiles = @(
    "go_stims-16515889668011746.csv",
    "GNG_stims-1651588966815421.csv",
    "go_stims-1651588966817431.csv",
    "GNG_stims-16515889668307655.csv"
)
$sheetGNG = $wb.Worksheets.Item("GNG_TO-16515889668327491")
for ($i = 0; $i -lt $gngFiles.Length; $i++) {
    $row = $i + 2
    $sheetGNG.Cells.Item($row, 2).Value = $gngFiles[$i]
}

# --- RS task order (rows 2-3, eyes closed / eyes open swapped) ---
$sheetRS = $wb.Worksheets.Item("RS_TO-16515889668348305")
$sheetRS.Cells.Item(2, 2).Value = "eyes closed"
$sheetRS.Cells.Item(3, 2).Value = "eyes open"

# --- vSAT task order (rows 2-5, column A index 0-3, column B file name) ---
$vsatFiles = @(
    "vSAT_stims-1651588966877771.csv",
    "SAT_stims-16515889668621376.csv",
    "vSAT_stims-16515889668928313.csv",
    "SAT_stims-1651588966839484.csv"
)
$sheetvSAT = $wb.Worksheets.Item("vSAT_TO-1651588966908666")
for ($i = 0; $i -lt $vsatFiles.Length; $i++) {
    $row = $i + 2
    $sheetvSAT.Cells.Item($row, 2).Value = $vsatFiles[$i]
}
